$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 16320
$ws.Range("I13").Value = 16250
$ws.Range("K13").Value = 16250
$ws.Range("M13").Value = -16081
$ws.Range("H28").Value = 445.17648
$ws.Range("I28").Value = 418.58334
$ws.Range("J28").Value = 509
$ws.Range("K28").Value = 418.58334
$ws.Range("L28").Value = 509
$ws.Range("M28").Value = 66.41665999999998
$ws.Range("N28").Value = -1479
$ws.Range("H138").Value = 2742668.5
$ws.Range("J138").Value = 3594.82
$ws.Range("L138").Value = 10784.46
$ws.Range("N138").Value = -21064.46

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10192.46
$ws.Range("I32").Value = 10719.237
$ws.Range("J32").Value = 2422.5
$ws.Range("K32").Value = 10719.237
$ws.Range("L32").Value = 2422.5
$ws.Range("M32").Value = -10432.237
$ws.Range("N32").Value = -2996.5
$ws.Range("H74").Value = 1138.8485
$ws.Range("I74").Value = 977.95
$ws.Range("K74").Value = 977.95
$ws.Range("M74").Value = -103.95
$ws.Range("H77").Value = 1138.8485
$ws.Range("I77").Value = 977.95
$ws.Range("K77").Value = 4889.75
$ws.Range("M77").Value = -521.75
$ws.Range("H102").Value = 1480.5
$ws.Range("I102").Value = 1475.625
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1475.625
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 146.375
$ws.Range("N102").Value = -4744
$ws.Range("H123").Value = 28428
$ws.Range("J123").Value = 28428
$ws.Range("L123").Value = 28428
$ws.Range("N123").Value = -38228
$ws.Range("H135").Value = 40655.75
$ws.Range("I135").Value = 50000
$ws.Range("J135").Value = 39320.855
$ws.Range("K135").Value = 50000
$ws.Range("L135").Value = 39320.855
$ws.Range("M135").Value = -44930
$ws.Range("N135").Value = -49460.855

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3402.6858
$ws.Range("I105").Value = 3306.7856
$ws.Range("J105").Value = 3466.6191
$ws.Range("K105").Value = 3306.7856
$ws.Range("L105").Value = 3466.6191
$ws.Range("M105").Value = -1559.7856
$ws.Range("N105").Value = -6960.6191
$ws.Range("H107").Value = 31460.234
$ws.Range("I107").Value = 37894.57
$ws.Range("J107").Value = 1433.3334
$ws.Range("K107").Value = 37894.57
$ws.Range("L107").Value = 1433.3334
$ws.Range("M107").Value = -35974.57
$ws.Range("N107").Value = -5273.3334
$ws.Range("H135").Value = 57017.5
$ws.Range("J135").Value = 57017.5
$ws.Range("L135").Value = 57017.5
$ws.Range("N135").Value = -67157.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1404.9
$ws.Range("I58").Value = 1373.8422
$ws.Range("J58").Value = 1995
$ws.Range("K58").Value = 1373.8422
$ws.Range("L58").Value = 1995
$ws.Range("M58").Value = -1170.8422
$ws.Range("N58").Value = -2401
$ws.Range("H86").Value = 4407.727
$ws.Range("I86").Value = 4164.1665
$ws.Range("J86").Value = 4700
$ws.Range("K86").Value = 4164.1665
$ws.Range("L86").Value = 4700
$ws.Range("M86").Value = -3041.1665
$ws.Range("N86").Value = -6946
$ws.Range("H89").Value = 4407.727
$ws.Range("I89").Value = 4164.1665
$ws.Range("J89").Value = 4700
$ws.Range("K89").Value = 20820.8325
$ws.Range("L89").Value = 23500
$ws.Range("M89").Value = -15204.8325
$ws.Range("N89").Value = -34732
$ws.Range("H105").Value = 1012.5
$ws.Range("I105").Value = 1000
$ws.Range("J105").Value = 1100
$ws.Range("K105").Value = 1000
$ws.Range("L105").Value = 1100
$ws.Range("M105").Value = 747
$ws.Range("N105").Value = -4594
$ws.Range("H107").Value = 432.16666
$ws.Range("I107").Value = 534.375
$ws.Range("J107").Value = 227.75
$ws.Range("K107").Value = 534.375
$ws.Range("L107").Value = 227.75
$ws.Range("M107").Value = 1385.625
$ws.Range("N107").Value = -4067.75
$ws.Range("H136").Value = 1404.9
$ws.Range("I136").Value = 1373.8422
$ws.Range("J136").Value = 1995
$ws.Range("K136").Value = 4121.5266
$ws.Range("L136").Value = 5985
$ws.Range("M136").Value = -1571.5266
$ws.Range("N136").Value = -11085

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null
$ws.Range("H103").Value = 26997.5
$ws.Range("J103").Value = 26997.5
$ws.Range("L103").Value = 26997.5
$ws.Range("N103").Value = -29341.5
$ws.Range("H113").Value = 1145.3334
$ws.Range("I113").Value = 782.55
$ws.Range("J113").Value = 1870.9
$ws.Range("K113").Value = 782.55
$ws.Range("L113").Value = 1870.9
$ws.Range("M113").Value = 1387.45
$ws.Range("N113").Value = -6210.9
$ws.Range("H122").Value = 3030.6316
$ws.Range("J122").Value = 2541.6
$ws.Range("L122").Value = 7624.799999999999
$ws.Range("N122").Value = -12524.8
$ws.Range("H123").Value = 12603.333
$ws.Range("J123").Value = 12603.333
$ws.Range("L123").Value = 12603.333
$ws.Range("N123").Value = -17503.333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5177.25
$ws.Range("I7").Value = 6352
$ws.Range("K7").Value = 6352
$ws.Range("M7").Value = -6240
$ws.Range("H61").Value = 18408.166
$ws.Range("I61").Value = 29485.428
$ws.Range("J61").Value = 2900
$ws.Range("K61").Value = 29485.428
$ws.Range("L61").Value = 2900
$ws.Range("M61").Value = -29283.428
$ws.Range("N61").Value = -3304
$ws.Range("H68").Value = 2036.625
$ws.Range("J68").Value = 2264.3333
$ws.Range("L68").Value = 2264.3333
$ws.Range("N68").Value = -3762.3333
$ws.Range("H71").Value = 2036.625
$ws.Range("J71").Value = 2264.3333
$ws.Range("L71").Value = 11321.6665
$ws.Range("N71").Value = -18809.6665
$ws.Range("H98").Value = 24751.666
$ws.Range("J98").Value = 24751.666
$ws.Range("L98").Value = 24751.666
$ws.Range("N98").Value = -30741.666
$ws.Range("H113").Value = 18408.166
$ws.Range("I113").Value = 29485.428
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 29485.428
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = -27315.428
$ws.Range("N113").Value = -7240
$ws.Range("H126").Value = 5177.25
$ws.Range("I126").Value = 6352
$ws.Range("K126").Value = 19056
$ws.Range("M126").Value = -16586
$ws.Range("H132").Value = 6359.2144
$ws.Range("I132").Value = 6153
$ws.Range("J132").Value = 6874.75
$ws.Range("K132").Value = 18459
$ws.Range("L132").Value = 20624.25
$ws.Range("M132").Value = -15929
$ws.Range("N132").Value = -25684.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1201.3529
$ws.Range("I113").Value = 954.9231
$ws.Range("J113").Value = 2002.25
$ws.Range("K113").Value = 2864.7693
$ws.Range("L113").Value = 6006.75
$ws.Range("M113").Value = -694.7692999999999
$ws.Range("N113").Value = -10346.75
$ws.Range("H122").Value = 27780690
$ws.Range("I122").Value = 31252900
$ws.Range("K122").Value = 93758700
$ws.Range("M122").Value = -93756250
$ws.Range("H126").Value = 14726.692
$ws.Range("I126").Value = 20827.445
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 62482.335
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -60012.335
$ws.Range("N126").Value = -7940

Write-Host "Applied all updates"